$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'" + '59.548.19'
$ws.Cells.Item(2, 5).Value = '  +0.59%  '

# Row 3
$ws.Cells.Item(3, 4).Value = "'" + '2.648.58'
$ws.Cells.Item(3, 5).Value = '  +0.23%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'" + '530.35'
$ws.Cells.Item(5, 5).Value = '  +1.55%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'" + '145.91'
$ws.Cells.Item(6, 5).Value = '  +0.35%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'" + '0.997'
$ws.Cells.Item(7, 5).Value = '  -0.25%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'" + '0.571'
$ws.Cells.Item(8, 5).Value = '  -0.55%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'" + '6.68'
$ws.Cells.Item(9, 5).Value = '  -3.34%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'" + '0.105'
$ws.Cells.Item(10, 5).Value = '  +2.09%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'" + '0.338'
$ws.Cells.Item(11, 5).Value = '  +0.90%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.61%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'" + '3.109.39'
$ws.Cells.Item(13, 5).Value = '  -0.03%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'" + '59.464.38'
$ws.Cells.Item(14, 5).Value = '  +0.36%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'" + '20.91'
$ws.Cells.Item(15, 5).Value = '  -1.17%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = "'" + '2.727.17'
$ws.Cells.Item(16, 5).Value = '  +2.33%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).Value = "'" + '0.0000137'
$ws.Cells.Item(17, 5).Value = '  +0.95%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'" + '343.34'
$ws.Cells.Item(18, 5).Value = '  +0.47%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'" + '4.47'
$ws.Cells.Item(19, 5).Value = '  +0.82%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'" + '10.66'
$ws.Cells.Item(20, 5).Value = '  +2.95%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'" + '6.41'
$ws.Cells.Item(21, 5).Value = '  +1.76%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'" + '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.08%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'" + '65.83'
$ws.Cells.Item(23, 5).Value = '  +3.34%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'" + '0.419'
$ws.Cells.Item(24, 5).Value = '  +1.63%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.20%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'" + '0.998'
$ws.Cells.Item(26, 5).Value = '  -0.23%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'" + '7.24'
$ws.Cells.Item(27, 5).Value = '  +1.74%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'" + '0.0₃0806'
$ws.Cells.Item(28, 5).Value = '  +0.34%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.05%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'" + '6.42'
$ws.Cells.Item(30, 5).Value = '  -3.74%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +1.86%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'" + '19.05'
$ws.Cells.Item(32, 5).Value = '  +1.45%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'" + '150.99'
$ws.Cells.Item(33, 5).Value = '  +1.11%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'" + '4.19'
$ws.Cells.Item(34, 5).Value = '  +0.08%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'" + '1.20'
$ws.Cells.Item(35, 5).Value = '  +0.42%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'" + '0.873'
$ws.Cells.Item(36, 5).Value = '  -2.10%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'" + '0.868'
$ws.Cells.Item(37, 5).Value = '  -1.68%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).Value = "'" + '1.50'
$ws.Cells.Item(38, 5).Value = '  +0.48%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'OKB'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(39, 4).Value = "'" + '36.56'
$ws.Cells.Item(39, 5).Value = '  -0.65%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'" + '3.66'
$ws.Cells.Item(40, 5).Value = '  +1.88%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.40%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Mantle'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(42, 4).Value = "'" + '0.605'
$ws.Cells.Item(42, 5).Value = '  -3.57%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Stellar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(43, 4).Value = "'" + '0.0978'
$ws.Cells.Item(43, 5).Value = '  +0.09%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'" + '270.96'
$ws.Cells.Item(44, 5).Value = '  -1.66%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'" + '19.49'
$ws.Cells.Item(45, 5).Value = '  -1.60%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +1.41%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'" + '0.0538'
$ws.Cells.Item(47, 5).Value = '  -0.35%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'" + '2.040.35'
$ws.Cells.Item(48, 5).Value = '  -0.50%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'" + '4.78'
$ws.Cells.Item(49, 5).Value = '  -0.03%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'" + '0.0230'
$ws.Cells.Item(50, 5).Value = '  +0.63%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'" + '18.88'
$ws.Cells.Item(51, 5).Value = '  -0.49%  '
